$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "statut" column (A) used black/green square emojis with labels "noir"/"vert".
# Refactor the synthetic array: swap the black-square status for a blue-book emoji
# with label "bleu" (green stays "vert" but its emoji changes too).
$ws.Range("A2").Value = "📘"
$ws.Range("B2").Value = "bleu"

$ws.Range("A3").Value = "📗"
$ws.Range("A4").Value = "📗"
